$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.625
$ws.Range("C2").Value = 0.3103448275862069
$ws.Range("D2").Value = 0.4147465437788019
$ws.Range("E2").Value = 145

# Row 3
$ws.Range("B3").Value = 0.8208955223880597
$ws.Range("C3").Value = 0.7432432432432432
$ws.Range("D3").Value = 0.7801418439716312
$ws.Range("E3").Value = 148

# Row 4
$ws.Range("B4").Value = 0.9537037037037037
$ws.Range("C4").Value = 0.6821192052980133
$ws.Range("D4").Value = 0.7953667953667954
$ws.Range("E4").Value = 151

# Row 5
$ws.Range("B5").Value = 0.4545454545454545
$ws.Range("C5").Value = 0.8333333333333334
$ws.Range("D5").Value = 0.5882352941176471
$ws.Range("E5").Value = 156

# Row 6 (accuracy row - all cells hold the same value)
$ws.Range("B6").Value = 0.6466666666666666
$ws.Range("C6").Value = 0.6466666666666666
$ws.Range("D6").Value = 0.6466666666666666
$ws.Range("E6").Value = 0.6466666666666666

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.7135361701593045
$ws.Range("C7").Value = 0.6422601523651992
$ws.Range("D7").Value = 0.6446226193087189

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.7117264791363049
$ws.Range("C8").Value = 0.6466666666666666
$ws.Range("D8").Value = 0.6457738895641112
